$d = $word.ActiveDocument

# --- Change 1: insert two new schedule lines after "-menu  18.11.2016"
# and before "System komentowania/wyszukiwania/oceniania" ---

$menuPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*-menu*18.11.2016*") {
        $menuPara = $p
        break
    }
}

$xmlTemplate = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/><w:t>{0}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$menuPara.Range.InsertParagraphAfter()
$newPara1 = $menuPara.Next()
$xml1 = $xmlTemplate -f "-layout pozostałych stron 25.11.2016"
[void]$newPara1.Range.InsertXML($xml1)

$newPara1 = $menuPara.Next()
$newPara1.Range.InsertParagraphAfter()
$newPara2 = $newPara1.Next()
$xml2 = $xmlTemplate -f "- dodawanie zdjęć  25.11.2016"
[void]$newPara2.Range.InsertXML($xml2)

# --- Change 2: drop the trailing blank paragraph and replace the
# "usuwanie z customowym komunikatem..." note with a single space ---

$notePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*usuwanie z customowym komunikatem*") {
        $notePara = $p
        break
    }
}

$blankPara = $notePara.Previous()
if ($blankPara.Range.Text.Trim().Length -eq 0) {
    $blankPara.Range.Delete()
}

# Re-resolve the paragraph: the earlier delete invalidates the cached
# range on $notePara, so look it up again before editing its text.
$notePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*usuwanie z customowym komunikatem*") {
        $notePara = $p
        break
    }
}
$notePara.Range.Text = " "
